$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "33.3%"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.0"

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "33.3%"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "5.7"
